$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 203, pushing the existing data
# (old rows 203-285) down to rows 205-287. Excel's row-insert copies
# formatting (e.g. the date style on column D) from the row above.
$ws.Rows("203:204").Insert()

# New row 203: "Primera" quality entry for 2022-12-22, Provincia de Curicó
$ws.Range("A203").Value2 = 3
$ws.Range("B203").Value2 = "Femacal de La Calera"
$ws.Range("C203").Value2 = "Coquimbo"
$ws.Range("D203").Value2 = 44917
$ws.Range("E203").Value2 = 5
$ws.Range("F203").Value2 = "Fruta"
$ws.Range("G203").Value2 = 100101
$ws.Range("H203").Value2 = "Berries"
$ws.Range("I203").Value2 = 100101001
$ws.Range("J203").Value2 = "Arándano (blue)"
$ws.Range("K203").Value2 = "Sin especificar"
$ws.Range("L203").Value2 = "Primera"
$ws.Range("M203").Value2 = 76
$ws.Range("N203").Value2 = 4600
$ws.Range("O203").Value2 = 4800
$ws.Range("P203").Value2 = 4705
$ws.Range("Q203").Value2 = "$/bandeja 2 kilos"
$ws.Range("R203").Value2 = "Provincia de Curicó"
$ws.Range("S203").Value2 = 2352
$ws.Range("T203").Value2 = 2

# New row 204: "Segunda" quality entry for 2022-12-22, Provincia de Curicó
$ws.Range("A204").Value2 = 3
$ws.Range("B204").Value2 = "Femacal de La Calera"
$ws.Range("C204").Value2 = "Coquimbo"
$ws.Range("D204").Value2 = 44917
$ws.Range("E204").Value2 = 5
$ws.Range("F204").Value2 = "Fruta"
$ws.Range("G204").Value2 = 100101
$ws.Range("H204").Value2 = "Berries"
$ws.Range("I204").Value2 = 100101001
$ws.Range("J204").Value2 = "Arándano (blue)"
$ws.Range("K204").Value2 = "Sin especificar"
$ws.Range("L204").Value2 = "Segunda"
$ws.Range("M204").Value2 = 90
$ws.Range("N204").Value2 = 3000
$ws.Range("O204").Value2 = 4500
$ws.Range("P204").Value2 = 3667
$ws.Range("Q204").Value2 = "$/bandeja 2 kilos"
$ws.Range("R204").Value2 = "Provincia de Curicó"
$ws.Range("S204").Value2 = 1834
$ws.Range("T204").Value2 = 2
